$wb = $excel.ActiveWorkbook

# The second sheet (index 2) is "Лист3" ( tabSelected in the original file ).
$ws3 = $wb.Worksheets.Item(2)

# Duplicate it right after itself -> becomes the new 3rd sheet ("Лист3 (2)"),
# matches the diff's new <sheet name="Лист4" .../> with sheetId 4 / rId3.
[void]$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(3)
$ws4.Name = "Лист4"

# --- Update the new "Лист4" sheet: Function 2 -> Function 3 block ---
$ws4.Range("B18").Value = "Function 3"
$ws4.Range("B19").Value = "Function 3"
$ws4.Range("B20").Value = "Function 3"
$ws4.Range("B21").Value = "Function 3"
$ws4.Range("B23").Value = "Function 3"
$ws4.Range("B24").Value = "Function 3"
$ws4.Range("B25").Value = "Function 3"
$ws4.Range("B26").Value = "Function 3"

$ws4.Range("E9").Value = "3;"
$ws4.Range("E10").Value = "                                                 +;"
$ws4.Range("E11").Value = "113;"
$ws4.Range("E12").Value = "4;"

$ws4.Range("M9").Value = -3
$ws4.Range("M10").Value = "*;"
$ws4.Range("M11").Value = "|;"
$ws4.Range("M12").Value = "*&%;"

$ws4.Range("D18").Value = "3. Print Team/s"
$ws4.Range("F18").Value = "3. Print Team/s"

# --- Update the original "Лист3" sheet row 18 text ---
$ws3.Range("D18").Value = "2. Print Teacher/s"
$ws3.Range("F18").Value = "2. Print Teacher/s"

# --- Selections: Лист3 keeps E28 selected, Лист4 becomes the active/selected tab with L20 selected ---
[void]$ws3.Activate()
[void]$ws3.Range("E28").Select()

[void]$ws4.Activate()
[void]$ws4.Range("L20").Select()
